$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Puerto Rico (AR) column for rows 97-99, matching the already-populated rows above (88-96)
$ws.Range("AR97").Value = 0
$ws.Range("AR98").Value = 0
$ws.Range("AR99").Value = 0

# Update recomputed values in row 119 (refreshed survey aggregate)
$ws.Range("B119").Value = 0.4992531
$ws.Range("D119").Value = 0.5054157
$ws.Range("G119").Value = 0.3227907
$ws.Range("R119").Value = 0.4132101
$ws.Range("S119").Value = 0.3934448
$ws.Range("U119").Value = 0.3723108
$ws.Range("V119").Value = 0.4391654
$ws.Range("AF119").Value = 0.3560332
$ws.Range("AG119").Value = 0.448961
$ws.Range("AJ119").Value = 0.4351192
$ws.Range("AP119").Value = 0.3279517
$ws.Range("AQ119").Value = 0.2745024
$ws.Range("AT119").Value = 0.357213
$ws.Range("AV119").Value = 0.363247
$ws.Range("AW119").Value = 0.331948
$ws.Range("BD119").Value = 0.3292094

# Append new survey rows 120-123 (29 May 2020 - 01 Jun 2020)
$ws.Range("A120").Value = "29 05 2020"
$ws.Range("B120").Value = 0.5324959
$ws.Range("C120").Value = 0.5408117
$ws.Range("D120").Value = 0.5350679
$ws.Range("F120").Value = 0.4963699
$ws.Range("G120").Value = 0.3210396
$ws.Range("H120").Value = 0.368475
$ws.Range("I120").Value = 0.3245705
$ws.Range("J120").Value = 0.6075444
$ws.Range("K120").Value = 0.3069542
$ws.Range("L120").Value = 0.3432742
$ws.Range("M120").Value = 0.3988088
$ws.Range("O120").Value = 0.1452381
$ws.Range("P120").Value = 0.5640774
$ws.Range("Q120").Value = 0.3247374
$ws.Range("R120").Value = 0.3919963
$ws.Range("S120").Value = 0.4142806
$ws.Range("T120").Value = 0.3373152
$ws.Range("U120").Value = 0.3695549
$ws.Range("V120").Value = 0.4008566
$ws.Range("W120").Value = 0.3446196
$ws.Range("X120").Value = 0.5314772
$ws.Range("Y120").Value = 0.2678815
$ws.Range("Z120").Value = 0.3130562
$ws.Range("AA120").Value = 0.4316825
$ws.Range("AB120").Value = 0.432488
$ws.Range("AD120").Value = 0.5949182
$ws.Range("AE120").Value = 0.3939767
$ws.Range("AF120").Value = 0.3498796
$ws.Range("AG120").Value = 0.3478986
$ws.Range("AH120").Value = 0.4421344
$ws.Range("AI120").Value = 0.2144979
$ws.Range("AJ120").Value = 0.4407877
$ws.Range("AK120").Value = 0.3526781
$ws.Range("AL120").Value = 0.2474397
$ws.Range("AM120").Value = 0.3457105
$ws.Range("AN120").Value = 0.360136
$ws.Range("AO120").Value = 0.3869271
$ws.Range("AP120").Value = 0.3532464
$ws.Range("AQ120").Value = 0.2865304
$ws.Range("AS120").Value = 0.4235876
$ws.Range("AT120").Value = 0.3496847
$ws.Range("AU120").Value = 0.535718
$ws.Range("AV120").Value = 0.364469
$ws.Range("AW120").Value = 0.3713251
$ws.Range("AX120").Value = 0.3217403
$ws.Range("AY120").Value = 0.3967036
$ws.Range("BA120").Value = 0.09051919999999999
$ws.Range("BB120").Value = 0.3265531
$ws.Range("BC120").Value = 0.2933674
$ws.Range("BD120").Value = 0.310627
$ws.Range("BE120").Value = 0.3492088

$ws.Range("A121").Value = "30 05 2020"
$ws.Range("B121").Value = 0.5047724
$ws.Range("C121").Value = 0.5242393
$ws.Range("D121").Value = 0.5960481
$ws.Range("F121").Value = 0.4692544
$ws.Range("G121").Value = 0.302697
$ws.Range("H121").Value = 0.3488608
$ws.Range("I121").Value = 0.3522608
$ws.Range("J121").Value = 0.5221980000000001
$ws.Range("K121").Value = 0.3075293
$ws.Range("L121").Value = 0.3424304
$ws.Range("M121").Value = 0.3850277
$ws.Range("O121").Value = 0.1453568
$ws.Range("P121").Value = 0.5872305
$ws.Range("Q121").Value = 0.4265771
$ws.Range("R121").Value = 0.3822531
$ws.Range("S121").Value = 0.4020259
$ws.Range("T121").Value = 0.330995
$ws.Range("U121").Value = 0.3695971
$ws.Range("V121").Value = 0.4590115
$ws.Range("W121").Value = 0.3384647
$ws.Range("X121").Value = 0.5344336
$ws.Range("Y121").Value = 0.2408678
$ws.Range("Z121").Value = 0.3203618
$ws.Range("AA121").Value = 0.3801123
$ws.Range("AB121").Value = 0.4461457
$ws.Range("AD121").Value = 0.5968539
$ws.Range("AE121").Value = 0.4369387
$ws.Range("AF121").Value = 0.3655755
$ws.Range("AG121").Value = 0.3767077
$ws.Range("AH121").Value = 0.4836817
$ws.Range("AI121").Value = 0.1938842
$ws.Range("AJ121").Value = 0.4122613
$ws.Range("AK121").Value = 0.3258516
$ws.Range("AL121").Value = 0.2748803
$ws.Range("AM121").Value = 0.3571308
$ws.Range("AN121").Value = 0.3354181
$ws.Range("AO121").Value = 0.4370689
$ws.Range("AP121").Value = 0.3593436
$ws.Range("AQ121").Value = 0.2840539
$ws.Range("AS121").Value = 0.4177374
$ws.Range("AT121").Value = 0.3393181
$ws.Range("AU121").Value = 0.4431948
$ws.Range("AV121").Value = 0.3900626
$ws.Range("AW121").Value = 0.3608041
$ws.Range("AX121").Value = 0.289303
$ws.Range("AY121").Value = 0.3896832
$ws.Range("BA121").Value = 0.0917873
$ws.Range("BB121").Value = 0.3014691
$ws.Range("BC121").Value = 0.2998637
$ws.Range("BD121").Value = 0.3496322
$ws.Range("BE121").Value = 0.3850579

$ws.Range("A122").Value = "31 05 2020"
$ws.Range("B122").Value = 0.4728497
$ws.Range("C122").Value = 0.509285
$ws.Range("D122").Value = 0.5498912
$ws.Range("F122").Value = 0.5341804999999999
$ws.Range("G122").Value = 0.3034673
$ws.Range("H122").Value = 0.3358322
$ws.Range("I122").Value = 0.3290508
$ws.Range("J122").Value = 0.4825046
$ws.Range("K122").Value = 0.361797
$ws.Range("L122").Value = 0.3555995
$ws.Range("M122").Value = 0.3554484
$ws.Range("O122").Value = 0.1480788
$ws.Range("P122").Value = 0.5314994
$ws.Range("Q122").Value = 0.4368779
$ws.Range("R122").Value = 0.3905919
$ws.Range("S122").Value = 0.4394632
$ws.Range("T122").Value = 0.341161
$ws.Range("U122").Value = 0.3498436
$ws.Range("V122").Value = 0.4467438
$ws.Range("W122").Value = 0.3496111
$ws.Range("X122").Value = 0.5275393
$ws.Range("Y122").Value = 0.228187
$ws.Range("Z122").Value = 0.2964295
$ws.Range("AA122").Value = 0.438177
$ws.Range("AB122").Value = 0.4163562
$ws.Range("AD122").Value = 0.6099749
$ws.Range("AE122").Value = 0.404754
$ws.Range("AF122").Value = 0.3428025
$ws.Range("AG122").Value = 0.4439807
$ws.Range("AH122").Value = 0.4408749
$ws.Range("AI122").Value = 0.2230886
$ws.Range("AJ122").Value = 0.3872506
$ws.Range("AK122").Value = 0.314305
$ws.Range("AL122").Value = 0.2907916
$ws.Range("AM122").Value = 0.3543519
$ws.Range("AN122").Value = 0.3508764
$ws.Range("AO122").Value = 0.4476105
$ws.Range("AP122").Value = 0.3237141
$ws.Range("AQ122").Value = 0.2874476
$ws.Range("AS122").Value = 0.3340212
$ws.Range("AT122").Value = 0.3199012
$ws.Range("AU122").Value = 0.4292347
$ws.Range("AV122").Value = 0.3254275
$ws.Range("AW122").Value = 0.3770988
$ws.Range("AX122").Value = 0.3158963
$ws.Range("AY122").Value = 0.387206
$ws.Range("BA122").Value = 0.15561
$ws.Range("BB122").Value = 0.3146316
$ws.Range("BC122").Value = 0.2959134
$ws.Range("BD122").Value = 0.4032489
$ws.Range("BE122").Value = 0.430363

$ws.Range("A123").Value = "01 06 2020"
